# Converts OOXML EMU (English Metric Units) to the points used by the
# PowerPoint COM object model (Shape.Left/Top/Width/Height), 1 pt = 12700 EMU.
# The host stores Left/Top/Width/Height as 32-bit floats and truncates
# (floor) pt*12700 back to EMU on save, so a plain division can land one
# EMU short after the float32 round-trip. Search nearby representable
# values and pick the one (closest to the true value) that truncates back
# to exactly the requested EMU amount.
function EmuToPt($emu) {
    $emuPerPt = 12700.0
    $base = [double]$emu / $emuPerPt
    $bestPt = $base
    $bestDist = 999999999.0
    for ($i = -200; $i -le 200; $i++) {
        $candD = $base + ($i * 0.00001)
        $candF = [float]$candD
        $back = [Math]::Floor([double]$candF * $emuPerPt)
        if ([int64]$back -eq $emu) {
            $dist = [Math]::Abs($candD - $base)
            if ($dist -lt $bestDist) {
                $bestDist = $dist
                $bestPt = $candD
            }
        }
    }
    return $bestPt
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape id=40 "TextBox 39" -> normalIWeightedString -> normalIWeightedStr, reposition
$sh = $s1.Shapes.Item(28)
$sh.Left = EmuToPt(3774707)
$sh.Top = EmuToPt(5915057)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 21).Text = "normalIWeightedStr"

# Shape id=41 "TextBox 40" -> normalItemString -> normalItemStr, resize width
$sh = $s1.Shapes.Item(29)
$sh.Width = EmuToPt(1872307)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 16).Text = "normalItemStr"

# Shape id=5 "TextBox 4" -> getWeightedString -> getWeightedStr, reposition + resize
$sh = $s1.Shapes.Item(37)
$sh.Left = EmuToPt(9919073)
$sh.Top = EmuToPt(5361252)
$sh.Width = EmuToPt(1329403)
$sh.TextFrame.TextRange.Text = "getWeightedStr"

# Shape id=71 "Rectangle 70" -> collapse "normalItemObj:Object" + "[] , " + "weightedItem" into "normalAndWeightedObj"
$sh = $s1.Shapes.Item(43)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 37).Text = "normalAndWeightedObj"

# Shape id=11 "Rectangle 10" -> weightedItem -> normalAndWeightedObj, reposition + resize
$sh = $s1.Shapes.Item(44)
$sh.Left = EmuToPt(7840961)
$sh.Top = EmuToPt(5039804)
$sh.Width = EmuToPt(1963423)
$sh.TextFrame.TextRange.Text = "normalAndWeightedObj"

# Shape id=72 "Rectangle 71" -> weightedItemString -> weightedItemStr, resize width
$sh = $s1.Shapes.Item(45)
$sh.Width = EmuToPt(1920013)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 18).Text = "weightedItemStr"

# Shape id=14 "Rectangle 13" -> normalItemObj -> normalAndWeightedObj, reposition + resize
$sh = $s1.Shapes.Item(46)
$sh.Left = EmuToPt(7832288)
$sh.Top = EmuToPt(5881024)
$sh.Width = EmuToPt(1963423)
$sh.TextFrame.TextRange.Text = "normalAndWeightedObj"

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape id=7 "TextBox 6": merge "normalItemObj," / "weightedItemobj: {" paragraphs
# into a single "normalAndWeightedObj : {" paragraph.
$sh = $s2.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(2, 1).Delete()
$tr.Characters(1, 13).Text = "normalAndWeightedObj"
$tr.Characters(21, 1).Text = " : {"
$sh.Width = EmuToPt(2147767)
